# Update the LAES course list to separate requirements into distinct
# columns: Prerequisites (existing, column C), Corequisites, Concurrent,
# Recommended (new), and Terms Typically Offered (shifted from D to G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "Terms Typically Offered"
# column (D), pushing it to G and shifting the sheet dimension to A1:G14.
$ws.Columns("D:F").Insert()

# New header labels for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Populate the new columns with "NA" for every existing data row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}
